# chore: update Sheets via scheduled runner
# Refreshes market-price / profit columns (H:N) for a handful of leve rows
# across the crafting job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

# row 38
$ws.Range("H38").Value = 4201.6
$ws.Range("I38").Value = 2008
$ws.Range("J38").Value = 4750
$ws.Range("K38").Value = 6024
$ws.Range("L38").Value = 14250
$ws.Range("M38").Value = -5652
$ws.Range("N38").Value = -14994

# row 42
$ws.Range("H42").Value = 308.83334
$ws.Range("I42").Value = 44.333332
$ws.Range("J42").Value = 573.3333
$ws.Range("K42").Value = 132.999996
$ws.Range("L42").Value = 1719.9999
$ws.Range("M42").Value = 97.00000399999999
$ws.Range("N42").Value = -2179.9999

# row 58
$ws.Range("H58").Value = 6683
$ws.Range("I58").Value = 305
$ws.Range("J58").Value = 9074.75
$ws.Range("K58").Value = 915
$ws.Range("L58").Value = 27224.25
$ws.Range("M58").Value = -765
$ws.Range("N58").Value = -27524.25

# row 129
$ws.Range("H129").Value = 971.3279
$ws.Range("I129").Value = 606.4761999999999
$ws.Range("J129").Value = 1162.875
$ws.Range("K129").Value = 1819.4286
$ws.Range("L129").Value = 3488.625
$ws.Range("M129").Value = 3180.5714
$ws.Range("N129").Value = -13488.625

# row 132
$ws.Range("H132").Value = 4017490
$ws.Range("I132").Value = 1113.0513
$ws.Range("J132").Value = 66672972
$ws.Range("K132").Value = 3339.1539
$ws.Range("L132").Value = 200018916
$ws.Range("M132").Value = -809.1539000000002
$ws.Range("N132").Value = -200023976

# row 138
$ws.Range("H138").Value = 4474.846
$ws.Range("I138").Value = 2291.6667
$ws.Range("J138").Value = 4994.651
$ws.Range("K138").Value = 6875.000100000001
$ws.Range("L138").Value = 14983.953
$ws.Range("M138").Value = -1735.000100000001
$ws.Range("N138").Value = -25263.953

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

# row 23
$ws.Range("H23").Value = 39029.375
$ws.Range("I23").Value = 40836.332
$ws.Range("J23").Value = 33608.5
$ws.Range("K23").Value = 40836.332
$ws.Range("L23").Value = 33608.5
$ws.Range("M23").Value = -40577.332
$ws.Range("N23").Value = -34126.5

# row 44 (N44 no longer applies - remove it)
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# row 55
$ws.Range("H55").Value = 20319.8
$ws.Range("J55").Value = 21899.75
$ws.Range("L55").Value = 21899.75
$ws.Range("N55").Value = -22529.75

# row 61
$ws.Range("H61").Value = 255099.47
$ws.Range("I61").Value = 5776.963
$ws.Range("J61").Value = 772923.1
$ws.Range("K61").Value = 5776.963
$ws.Range("L61").Value = 772923.1
$ws.Range("M61").Value = -5564.963
$ws.Range("N61").Value = -773347.1

# row 123 (new N123 now applies)
$ws.Range("H123").Value = 34429
$ws.Range("J123").Value = 34429
$ws.Range("L123").Value = 34429
$ws.Range("N123").Value = -44229

# row 132
$ws.Range("H132").Value = 3159.0222
$ws.Range("I132").Value = 2166.28
$ws.Range("J132").Value = 4399.95
$ws.Range("K132").Value = 6498.84
$ws.Range("L132").Value = 13199.85
$ws.Range("M132").Value = -3968.84
$ws.Range("N132").Value = -18259.85

# row 136
$ws.Range("H136").Value = 255099.47
$ws.Range("I136").Value = 5776.963
$ws.Range("J136").Value = 772923.1
$ws.Range("K136").Value = 17330.889
$ws.Range("L136").Value = 2318769.3
$ws.Range("M136").Value = -14780.889
$ws.Range("N136").Value = -2323869.3

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

# row 134
$ws.Range("H134").Value = 29970.795
$ws.Range("I134").Value = 4443.3423
$ws.Range("K134").Value = 13330.0269
$ws.Range("M134").Value = -10795.0269

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

# row 31
$ws.Range("H31").Value = 16671002
$ws.Range("I31").Value = 1465.579
$ws.Range("K31").Value = 1465.579
$ws.Range("M31").Value = -1170.579

# row 34
$ws.Range("H34").Value = 16671002
$ws.Range("I34").Value = 1465.579
$ws.Range("K34").Value = 1465.579
$ws.Range("M34").Value = -1263.579

# row 134
$ws.Range("H134").Value = 8806022
$ws.Range("I134").Value = 12348953
$ws.Range("K134").Value = 37046859
$ws.Range("M134").Value = -37044324

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

# row 113
$ws.Range("H113").Value = 156742.36
$ws.Range("I113").Value = 482.81818
$ws.Range("J113").Value = 238592.6
$ws.Range("K113").Value = 1448.45454
$ws.Range("L113").Value = 715777.8
$ws.Range("M113").Value = 721.54546
$ws.Range("N113").Value = -720117.8

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

# row 43
$ws.Range("H43").Value = 3475
$ws.Range("I43").Value = 1800
$ws.Range("J43").Value = 8500
$ws.Range("K43").Value = 1800
$ws.Range("L43").Value = 8500
$ws.Range("M43").Value = -1649
$ws.Range("N43").Value = -8802

# row 57
$ws.Range("H57").Value = 6044.4116
$ws.Range("J57").Value = 6044.4116
$ws.Range("L57").Value = 6044.4116
$ws.Range("N57").Value = -7684.4116

# row 80
$ws.Range("H80").Value = 9337.571
$ws.Range("I80").Value = 21684
$ws.Range("J80").Value = 2478.4443
$ws.Range("K80").Value = 21684
$ws.Range("L80").Value = 2478.4443
$ws.Range("M80").Value = -20686
$ws.Range("N80").Value = -4474.4443

# row 83
$ws.Range("H83").Value = 9337.571
$ws.Range("I83").Value = 21684
$ws.Range("J83").Value = 2478.4443
$ws.Range("K83").Value = 108420
$ws.Range("L83").Value = 12392.2215
$ws.Range("M83").Value = -103428
$ws.Range("N83").Value = -22376.2215

# row 132
$ws.Range("H132").Value = 6668599
$ws.Range("I132").Value = 7577362.5
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 22732087.5
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -22729557.5
$ws.Range("N132").Value = -18059

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

# row 74 (new M74 now applies)
$ws.Range("H74").Value = 20250
$ws.Range("I74").Value = 18000
$ws.Range("J74").Value = 22500
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 22500
$ws.Range("M74").Value = -17002
$ws.Range("N74").Value = -24496

# row 77 (new M77 now applies)
$ws.Range("H77").Value = 20250
$ws.Range("I77").Value = 18000
$ws.Range("J77").Value = 22500
$ws.Range("K77").Value = 54000
$ws.Range("L77").Value = 67500
$ws.Range("M77").Value = -49008
$ws.Range("N77").Value = -77484

# row 122
$ws.Range("H122").Value = 2649054.2
$ws.Range("I122").Value = 3109129
$ws.Range("J122").Value = 3623.75
$ws.Range("K122").Value = 9327387
$ws.Range("L122").Value = 10871.25
$ws.Range("M122").Value = -9324937
$ws.Range("N122").Value = -15771.25

# row 137 (new M137 now applies)
$ws.Range("H137").Value = 54138.168
$ws.Range("I137").Value = 50000
$ws.Range("J137").Value = 54965.8
$ws.Range("K137").Value = 50000
$ws.Range("L137").Value = 54965.8
$ws.Range("M137").Value = -44900
$ws.Range("N137").Value = -65165.8

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

# row 81
$ws.Range("H81").Value = 12564019
$ws.Range("J81").Value = 33501968
$ws.Range("L81").Value = 67003936
$ws.Range("N81").Value = -67006058

# row 84
$ws.Range("H84").Value = 12564019
$ws.Range("J84").Value = 33501968
$ws.Range("L84").Value = 335019680
$ws.Range("N84").Value = -335030288

# row 122
$ws.Range("H122").Value = 1213.6
$ws.Range("I122").Value = 922.6667
$ws.Range("K122").Value = 2768.0001
$ws.Range("M122").Value = -318.0001000000002

# row 136
$ws.Range("H136").Value = 3548618
$ws.Range("I136").Value = 2839.423
$ws.Range("J136").Value = 7938629.5
$ws.Range("K136").Value = 8518.269
$ws.Range("L136").Value = 23815888.5
$ws.Range("M136").Value = -5968.269
$ws.Range("N136").Value = -23820988.5
